$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Rename sheet1 ---
$ws1.Name = "POS_Login"

# --- sheet1 (POS_Login) content ---
$ws1.Range("A1").Value = "UserName"
$ws1.Range("B1").Value = "Password"
$ws1.Range("A2").Value = "admin"
$ws1.Range("B2").Value = "nopoint"
$ws1.Range("A3").Value = "admin"
$ws1.Range("A4").Value = "bhanu"
$ws1.Range("B4").Value = "pointofsale"
$ws1.Range("A5").Value = "bhanu"
$ws1.Range("A6").Value = "admin"
$ws1.Range("B6").Value = "POINTOFSALE"
$ws1.Range("B8").Value = "pointofsale"

# --- sheet1 styling: style 1 = thin border only ---
$ws1.Range("A1:B2").Borders.LineStyle = 1
$ws1.Range("A1:B2").Borders.ColorIndex = 64
$ws1.Range("A3").Borders.LineStyle = 1
$ws1.Range("A3").Borders.ColorIndex = 64
$ws1.Range("A4:B4").Borders.LineStyle = 1
$ws1.Range("A4:B4").Borders.ColorIndex = 64
$ws1.Range("A5").Borders.LineStyle = 1
$ws1.Range("A5").Borders.ColorIndex = 64
$ws1.Range("A6:B6").Borders.LineStyle = 1
$ws1.Range("A6:B6").Borders.ColorIndex = 64
$ws1.Range("B8").Borders.LineStyle = 1
$ws1.Range("B8").Borders.ColorIndex = 64

# --- sheet1 styling: style 2 = thin border + yellow fill (empty highlighted cells) ---
$ws1.Range("B3").Borders.LineStyle = 1
$ws1.Range("B3").Borders.ColorIndex = 64
$ws1.Range("B3").Interior.Color = 65535

$ws1.Range("B5").Borders.LineStyle = 1
$ws1.Range("B5").Borders.ColorIndex = 64
$ws1.Range("B5").Interior.Color = 65535

$ws1.Range("A7:B7").Borders.LineStyle = 1
$ws1.Range("A7:B7").Borders.ColorIndex = 64
$ws1.Range("A7:B7").Interior.Color = 65535

$ws1.Range("A8").Borders.LineStyle = 1
$ws1.Range("A8").Borders.ColorIndex = 64
$ws1.Range("A8").Interior.Color = 65535

# --- sheet1 column widths: auto-fit to the new (longer) content ---
$ws1.Columns.Item(1).AutoFit()
$ws1.Columns.Item(2).AutoFit()

# --- sheet2 (InvalidLogin) content: keep same text, shared-string indices shift automatically ---
$ws2.Range("A2").Value = "abcd123"
$ws2.Range("B2").Value = "xyz456"
$ws2.Range("A3").Value = "admin"
$ws2.Range("B3").Value = "damager"

# --- selections & active sheet/tab ---
$ws2.Range("C7").Select()
$ws1.Activate()
$ws1.Range("B5").Select()
